$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "GPIO"

$ws.Range("A1").Value = "GPIOx"
$ws.Range("B1").Value = "PINx"
$ws.Range("C1").Value = "IN/OUT"
$ws.Range("D1").Value = "PULL/DROP"
$ws.Range("E1").Value = "SET/RESET"

$ws.Range("A2").Value = "GPIOA"
$ws.Range("B2").Value = "PIN10"
$ws.Range("C2").Value = "OUT"
$ws.Range("D2").Value = "PULL"
$ws.Range("E2").Value = "SET"

$ws.Range("A2").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
